$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.792.34"
$ws.Range("E2").Value = "  -0.47%  "
$ws.Range("D3").Value = "2.343.17"
$ws.Range("E3").Value = "  -0.38%  "
$ws.Range("E4").Value = "  -0.20%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "238.94"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.09%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.663"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.98%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "72.63"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.68%  "
$ws.Range("E8").Value = "  -0.06%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.597"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.86%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.100"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.40%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "60.81"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +6.38%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "32.98"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.65%  "
$ws.Range("E13").Value = "  +0.12%  "
$ws.Range("E14").Value = "  -0.46%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "16.04"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.00%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.898"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.04%  "
$ws.Range("D17").Value = "2.343.84"
$ws.Range("E17").Value = "  -1.26%  "
$ws.Range("D18").Value = "43.747.81"
$ws.Range("E18").Value = "  -0.35%  "
$ws.Range("E19").Value = "  +0.08%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "77.63"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.24%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.46"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.26%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "251.39"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.75%  "
$ws.Range("E23").Value = "  +2.85%  "
$ws.Range("E24").Value = "  -0.04%  "
$ws.Range("E25").Value = "  -4.98%  "
$ws.Range("E26").Value = "  -0.13%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.35"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.81%  "
$ws.Range("E28").Value = "  +0.79%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "175.42"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.10%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "22.11"
$ws.Range("D30").Style = "Normal"
$ws.Range("E31").Value = "  +0.32%  "
$ws.Range("E32").Value = "  -2.16%  "
$ws.Range("E33").Value = "  -2.34%  "
$ws.Range("E34").Value = "  -4.62%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.30"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.58%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.73"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.03%  "
$ws.Range("B37").Value = "LidoDAOToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.38"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.05%  "
$ws.Range("B38").Value = "THORChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.38"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.76%  "
$ws.Range("B39").Value = "FTXToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.58"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +19.50%  "
$ws.Range("B40").Value = "VeChain"
$ws.Range("C40").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0270"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.65%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "65.31"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +13.72%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "19.66"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.65%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "9.02"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.58%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.105"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.90%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.197"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.70%  "
$ws.Range("E46").Value = "  +0.04%  "
$ws.Range("E47").Value = "  -1.94%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.40"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.03%  "
$ws.Range("E49").Value = "  -2.09%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "97.16"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.79%  "
$ws.Range("E51").Value = "  +1.84%  "
